# Add three new network-device rows (SW4, SW5, SW6) to the address table,
# mirroring the formatting of the existing SW1/SW2/SW3 rows (25-27),
# then move the selection to I18 as in the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 28; A = "SW4"; B = "Vlan1"; C = "192.168.27.7";   D = "255.255.255.128" },
    @{ Row = 29; A = "SW5"; B = "Vlan1"; C = "192.168.28.2";   D = "255.255.255.192" },
    @{ Row = 30; A = "SW6"; B = "Vlan1"; C = "192.168.28.130"; D = "255.255.255.192" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Value = $r.A
    $cellA.HorizontalAlignment = -4108  # xlCenter (matches A25:A27 style)

    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
}

# Match the saved selection / scroll state recorded in the workbook.
$ws.Range("I18").Select() | Out-Null
